$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes --------------------------------------------------
# 1) Drop the old "Constant"/"r2_adj" row (row 5) entirely.
$ws.Rows(5).Delete()

# 2) Insert a new column between the existing "C" and "$\pi$" columns for
#    the new "A" regressor; this shifts the old column C -> D while
#    re-using the existing border/bold style (no new style is minted).
$ws.Columns("C").Insert()

# 3) Make room for a new "A Lag" row between "C Lag" and "$\pi$ Lag" by
#    shifting the "$\pi$ Lag" row (row 3) down into row 4 (which held the
#    now-removed "Constant" row) via copy/paste so styles are re-used
#    instead of synthesizing new ones via Insert().
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial(-4104)

# Re-apply the existing bold/border label style to A3 (it still carries
# the stale "$\pi$ Lag" value/style from the row above the shift; copy the
# format from A2, which already uses that exact style index).
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# --- Content ---------------------------------------------------------------
# Helper: write a plain piece of text into a cell without ever touching its
# NumberFormat/Font/Border (any of those mint a brand-new style entry in
# this runtime). Values that Excel wouldn't auto-parse as a number can be
# set directly via .Value. Purely-numeric-looking text (e.g. "1.004")
# would otherwise be silently stored as a Number, so instead we compute it
# with TEXT() and then Paste-Special "Values only" over the formula - that
# converts the cell to a literal text value while leaving style alone.

# Header row
$ws.Range("B1").Value = "C"
$ws.Range("C1").Value = "A"
$ws.Range("D1").Value = "$\pi$"

# Row 2 - C Lag
$ws.Range("A2").Value = "C Lag"
$ws.Range("B2").Value = "-0.698***"
$ws.Range("C2").Value = "-0.052***"
$ws.Range("D2").Formula = '=TEXT(-0.013,"0.000")'

# Row 3 - A Lag (new)
$ws.Range("A3").Value = "A Lag"
$ws.Range("B3").Formula = '=TEXT(1.004,"0.000")'
$ws.Range("C3").Value = "-0.412***"
$ws.Range("D3").Formula = '=TEXT(0.202,"0.000")'

# Row 4 - $\pi$ Lag
$ws.Range("A4").Value = "$\pi$ Lag"
$ws.Range("B4").Value = "1.118**"
$ws.Range("C4").Formula = '=TEXT(-0.071,"0.000")'
$ws.Range("D4").Value = "-0.636***"

# Freeze the TEXT()-formula cells down to their literal string results
# (one at a time - this runtime's PasteSpecial does not reliably apply to
# multi-area union ranges).
foreach ($addr in @("D2", "B3", "D3", "C4")) {
    $cell = $ws.Range($addr)
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
